# Component Analysis filtering: the rolling evaluation window advanced by one
# period, so each existing data row (B:G) shifts down into the next row and a
# new row of figures is written at the top of the data block (row 2). The
# oldest row (which was at row 11) drops out of the fixed-size window.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataCol = 2   # column B
$lastDataCol = 7    # column G
$firstDataRow = 2
$lastDataRow = 11

# Shift existing B:G values down by one row (start from the bottom so we
# don't clobber a row before it has been copied downward).
for ($r = $lastDataRow; $r -gt $firstDataRow; $r--) {
    $srcRow = $r - 1
    for ($c = $firstDataCol; $c -le $lastDataCol; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $ws.Cells.Item($srcRow, $c).Value2
    }
}

# New values for the top row of the data block (row 2).
$newRow2 = @(0.1369420173923726, 0.9020021672123393, 4.56720535997291, 2.137102093951739, 2.15627664282098, 46)
for ($c = $firstDataCol; $c -le $lastDataCol; $c++) {
    $ws.Cells.Item($firstDataRow, $c).Value2 = $newRow2[$c - $firstDataCol]
}
